$d = $word.ActiveDocument

# Paragraph 1 (Title): merge the per-word runs into a single run.
$p1 = $d.Paragraphs.Item(1)
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Answers: Introduction to simultaneous equations</w:t></w:r></w:p>'
$null = $p1.Range.InsertXML($xml1)

# Paragraph 2 (Author): merge the per-word runs into a single run.
$p2 = $d.Paragraphs.Item(2)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Ollie Brooke</w:t></w:r></w:p>'
$null = $p2.Range.InsertXML($xml2)

# Paragraph 4 (Abstract): merge the per-word runs into a single run.
$p4 = $d.Paragraphs.Item(4)
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Abstract"/></w:pPr><w:r><w:t xml:space="preserve">Answers to questions relating to the guide on introduction to simultaneous equations.</w:t></w:r></w:p>'
$null = $p4.Range.InsertXML($xml4)
